$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" (Changed) date in column C from 2023-09-08 (45177)
# to 2023-09-09 (45178) for all data rows (rows 2-46).
for ($r = 2; $r -le 46; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value = 45178
    }
}
